$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.012.54"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "1.955.13"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.02"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4879"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2959"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06815"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.18"
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.44"
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("D12").Value = "1.954.00"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.453"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7038"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.00"
$ws.Range("E16").Value = "  -3.24%  "
$ws.Range("D17").Value = "31.029.98"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.18"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007685"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "2.212.98"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.501"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.494"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.790"
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.48"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.00"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.212"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1059"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.404"
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.461"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04923"
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7638"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.170"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.727"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.701"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.511"
$ws.Range("E40").Value = "  +4.61%  "
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "75.52"
$ws.Range("E42").Value = "  +7.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8927"
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4471"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "109.39"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.212"
$ws.Range("E46").Value = "  +11.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9996"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "991.98"
$ws.Range("E48").Value = "  +8.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1258"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.351"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.73"
$ws.Range("E51").Value = "  -0.75%  "
